$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos price (D) and 1h-volume-change (E) refresh.
# Numeric-looking D values are forced to Text (matching the source
# inlineStr cells) by staging NumberFormat "@" before the write, then
# resetting the style back to Normal so no extra formatting sticks.

$ws.Range("D2").Value = '26.439.01'
$ws.Range("E2").Value = '  -2.66%  '

$ws.Range("D3").Value = '1.774.05'
$ws.Range("E3").Value = '  -1.60%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("E5").Value = '  +0.07%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '306.29'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.11%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4274'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.75%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3626'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.59%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07142'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.57%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8391'
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.53'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.34%  '

$ws.Range("D12").Value = '1.807.30'
$ws.Range("E12").Value = '  -3.25%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.442'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.51%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.241'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.16%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06894'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.14%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.009'
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '78.71'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.63%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008685'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.59%  '

$ws.Range("E19").Value = '  -0.06%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.93'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.51%  '

$ws.Range("D21").Value = '26.450.81'
$ws.Range("E21").Value = '  -2.72%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.102'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.70%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.06'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.07%  '

$ws.Range("D24").Value = '2.022.86'
$ws.Range("E24").Value = '  -2.72%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.48'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.55%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.857'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.08%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.99'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.89%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.048'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.33%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '113.83'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.82%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.773'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.35%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08888'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.18%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7233'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.71%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.114'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.66%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.306'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.06%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.003'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.05%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.747'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.70%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.100'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.98%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05127'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.32%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01886'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.33%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.1610'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.22%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4902'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.95%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.610'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.17%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.326'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.93%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.956'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.18%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '104.62'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.03%  '

$ws.Range("E46").Value = '  +0.17%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.12'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.53%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.628'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.65%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06182'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.24%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4464'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.71%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.708'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.37%  '
